$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 109.333336
$ws.Range("I12").Value = 109.333336
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 109.333336
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 60.666664
$ws.Range("N12").Value = $null
$ws.Range("H17").Value = 387.6279
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 387.6279
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1162.8837
$ws.Range("N17").Value = -1498.8837
$ws.Range("H74").Value = 4112.5386
$ws.Range("I74").Value = 3494.1428
$ws.Range("J74").Value = 4834
$ws.Range("K74").Value = 3494.1428
$ws.Range("L74").Value = 4834
$ws.Range("M74").Value = -2558.1428
$ws.Range("N74").Value = -6706
$ws.Range("H77").Value = 4112.5386
$ws.Range("I77").Value = 3494.1428
$ws.Range("J77").Value = 4834
$ws.Range("K77").Value = 17470.714
$ws.Range("L77").Value = 24170
$ws.Range("M77").Value = -12790.714
$ws.Range("N77").Value = -33530
$ws.Range("H87").Value = 11561.289
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 11561.289
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 11561.289
$ws.Range("N87").Value = -14057.289
$ws.Range("H90").Value = 11561.289
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 11561.289
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 34683.867
$ws.Range("N90").Value = -47163.867
$ws.Range("H100").Value = 11496063
$ws.Range("I100").Value = 15874025
$ws.Range("J100").Value = 3912.5
$ws.Range("K100").Value = 15874025
$ws.Range("L100").Value = 3912.5
$ws.Range("M100").Value = -15873484
$ws.Range("N100").Value = -4994.5
$ws.Range("H113").Value = 2454.6072
$ws.Range("I113").Value = 2419.5217
$ws.Range("J113").Value = 2616
$ws.Range("K113").Value = 2419.5217
$ws.Range("L113").Value = 2616
$ws.Range("M113").Value = 834.4783000000002
$ws.Range("N113").Value = -9124
$ws.Range("H132").Value = 11907170
$ws.Range("I132").Value = 11907170
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 35721510
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -35718980
$ws.Range("N132").Value = $null
$ws.Range("H140").Value = 67920
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 67920
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 67920
$ws.Range("N140").Value = -78280

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 338333
$ws.Range("I32").Value = 3177.7703
$ws.Range("J32").Value = 2405123.8
$ws.Range("K32").Value = 3177.7703
$ws.Range("L32").Value = 2405123.8
$ws.Range("M32").Value = -2890.7703
$ws.Range("N32").Value = -2405697.8
$ws.Range("H37").Value = 1300
$ws.Range("I37").Value = 1300
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1300
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -1027
$ws.Range("N37").Value = $null
$ws.Range("H61").Value = 4641.5454
$ws.Range("I61").Value = 4400
$ws.Range("J61").Value = 4931.4
$ws.Range("K61").Value = 4400
$ws.Range("L61").Value = 4931.4
$ws.Range("M61").Value = -4188
$ws.Range("N61").Value = -5355.4
$ws.Range("H63").Value = 3612.1428
$ws.Range("I63").Value = 3274.6155
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 3274.6155
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -2588.6155
$ws.Range("H66").Value = 3612.1428
$ws.Range("I66").Value = 3274.6155
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 16373.0775
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -12941.0775
$ws.Range("H132").Value = 12847855
$ws.Range("I132").Value = 19232244
$ws.Range("J132").Value = 79077.766
$ws.Range("K132").Value = 57696732
$ws.Range("L132").Value = 237233.298
$ws.Range("M132").Value = -57694202
$ws.Range("N132").Value = -242293.298
$ws.Range("H136").Value = 4641.5454
$ws.Range("I136").Value = 4400
$ws.Range("J136").Value = 4931.4
$ws.Range("K136").Value = 13200
$ws.Range("L136").Value = 14794.2
$ws.Range("M136").Value = -10650
$ws.Range("N136").Value = -19894.2

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1077744
$ws.Range("I31").Value = 2151947.2
$ws.Range("J31").Value = 3540.7097
$ws.Range("K31").Value = 2151947.2
$ws.Range("L31").Value = 3540.7097
$ws.Range("M31").Value = -2151652.2
$ws.Range("N31").Value = -4130.709699999999
$ws.Range("H34").Value = 1077744
$ws.Range("I34").Value = 2151947.2
$ws.Range("J34").Value = 3540.7097
$ws.Range("K34").Value = 2151947.2
$ws.Range("L34").Value = 3540.7097
$ws.Range("M34").Value = -2151745.2
$ws.Range("N34").Value = -3944.7097
$ws.Range("H58").Value = 1445
$ws.Range("I58").Value = 1410.1666
$ws.Range("J58").Value = 1549.5
$ws.Range("K58").Value = 1410.1666
$ws.Range("L58").Value = 1549.5
$ws.Range("M58").Value = -1207.1666
$ws.Range("N58").Value = -1955.5
$ws.Range("H136").Value = 1445
$ws.Range("I136").Value = 1410.1666
$ws.Range("J136").Value = 1549.5
$ws.Range("K136").Value = 4230.4998
$ws.Range("L136").Value = 4648.5
$ws.Range("M136").Value = -1680.4998
$ws.Range("N136").Value = -9748.5
$ws.Range("H138").Value = 49966.668
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 49966.668
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 49966.668
$ws.Range("N138").Value = -60246.668

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 385.77777
$ws.Range("I2").Value = 96.333336
$ws.Range("J2").Value = 443.66666
$ws.Range("K2").Value = 578.000016
$ws.Range("L2").Value = 2661.99996
$ws.Range("M2").Value = -465.000016
$ws.Range("N2").Value = -2887.99996
$ws.Range("H3").Value = 600
$ws.Range("I3").Value = 600
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1688
$ws.Range("H7").Value = 76.45
$ws.Range("I7").Value = 9
$ws.Range("J7").Value = 98.933334
$ws.Range("K7").Value = 27
$ws.Range("L7").Value = 296.800002
$ws.Range("M7").Value = 85
$ws.Range("N7").Value = -520.8000019999999
$ws.Range("H33").Value = 311.95
$ws.Range("I33").Value = 124.1875
$ws.Range("J33").Value = 1063
$ws.Range("K33").Value = 745.125
$ws.Range("L33").Value = 6378
$ws.Range("M33").Value = -462.125
$ws.Range("N33").Value = -6944
$ws.Range("H34").Value = 611.9286
$ws.Range("I34").Value = 172.5
$ws.Range("J34").Value = 787.7
$ws.Range("K34").Value = 517.5
$ws.Range("L34").Value = 2363.1
$ws.Range("M34").Value = -433.5
$ws.Range("N34").Value = -2531.1
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").Value = $null
$ws.Range("H51").Value = 1580.8334
$ws.Range("I51").Value = 500
$ws.Range("J51").Value = 2121.25
$ws.Range("K51").Value = 1500
$ws.Range("L51").Value = 6363.75
$ws.Range("M51").Value = -1040
$ws.Range("N51").Value = -7283.75
$ws.Range("H68").Value = 1234.079
$ws.Range("I68").Value = 773.3333
$ws.Range("J68").Value = 1534.5652
$ws.Range("K68").Value = 2319.9999
$ws.Range("L68").Value = 4603.6956
$ws.Range("M68").Value = -1508.9999
$ws.Range("N68").Value = -6225.6956
$ws.Range("H71").Value = 1234.079
$ws.Range("I71").Value = 773.3333
$ws.Range("J71").Value = 1534.5652
$ws.Range("K71").Value = 6959.9997
$ws.Range("L71").Value = 13811.0868
$ws.Range("M71").Value = -2903.9997
$ws.Range("N71").Value = -21923.0868
$ws.Range("H113").Value = 846.59753
$ws.Range("I113").Value = 589.94446
$ws.Range("J113").Value = 918.78125
$ws.Range("K113").Value = 1769.83338
$ws.Range("L113").Value = 2756.34375
$ws.Range("M113").Value = 400.16662
$ws.Range("N113").Value = -7096.34375
$ws.Range("H122").Value = 382.86667
$ws.Range("I122").Value = 326.08
$ws.Range("J122").Value = 666.8
$ws.Range("K122").Value = 2934.72
$ws.Range("L122").Value = 6001.2
$ws.Range("M122").Value = -484.7199999999998
$ws.Range("N122").Value = -10901.2
$ws.Range("H131").Value = 8929642
$ws.Range("I131").Value = 1177.1666
$ws.Range("J131").Value = 10001058
$ws.Range("K131").Value = 3531.4998
$ws.Range("L131").Value = 30003174
$ws.Range("M131").Value = 1508.5002
$ws.Range("N131").Value = -30013254
$ws.Range("H132").Value = 2980
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2980
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 26820
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -31880
$ws.Range("H140").Value = 7572.6
$ws.Range("I140").Value = 6465.8335
$ws.Range("J140").Value = 11999.667
$ws.Range("K140").Value = 19397.5005
$ws.Range("L140").Value = 35999.001
$ws.Range("M140").Value = -14217.5005
$ws.Range("N140").Value = -46359.001

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1433.8
$ws.Range("I113").Value = 1133.3334
$ws.Range("J113").Value = 1634.1111
$ws.Range("K113").Value = 1133.3334
$ws.Range("L113").Value = 1634.1111
$ws.Range("M113").Value = 1036.6666
$ws.Range("N113").Value = -5974.1111
$ws.Range("H138").Value = 58000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 58000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 58000
$ws.Range("N138").Value = -68280
$ws.Range("H140").Value = 50000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 50000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1673.2632
$ws.Range("I7").Value = 1492.6666
$ws.Range("J7").Value = 1982.8572
$ws.Range("K7").Value = 1492.6666
$ws.Range("L7").Value = 1982.8572
$ws.Range("M7").Value = -1380.6666
$ws.Range("N7").Value = -2206.8572
$ws.Range("H126").Value = 1673.2632
$ws.Range("I126").Value = 1492.6666
$ws.Range("J126").Value = 1982.8572
$ws.Range("K126").Value = 4477.9998
$ws.Range("L126").Value = 5948.571599999999
$ws.Range("M126").Value = -2007.9998
$ws.Range("N126").Value = -10888.5716

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 90286
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 90286
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 90286
$ws.Range("N46").Value = -90748
$ws.Range("H132").Value = 35717908
$ws.Range("I132").Value = 57693972
$ws.Range("J132").Value = 6808.4585
$ws.Range("K132").Value = 173081916
$ws.Range("L132").Value = 20425.3755
$ws.Range("M132").Value = -173079386
$ws.Range("N132").Value = -25485.3755
$ws.Range("H134").Value = 90286
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 90286
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 270858
$ws.Range("N134").Value = -275928
$ws.Range("H136").Value = 63219
$ws.Range("I136").Value = 111657.11
$ws.Range("J136").Value = 941.4286
$ws.Range("K136").Value = 334971.33
$ws.Range("L136").Value = 2824.2858
$ws.Range("M136").Value = -332421.33
$ws.Range("N136").Value = -7924.2858
